# Add a team record (Wins / Losses / Ties) to the data, per the commit:
# "Added team record to data" -- the W/L/T live on the same sheet as the
# player data, in three new trailing columns (AD, AE, AF).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new headers, matching the existing header style ---
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Reuse the exact header formatting (font/border/alignment) already used
# for the other header cells, e.g. AC1, instead of building a new style.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# PasteSpecial only touches formatting; make sure the text values still hold.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (2-62): the team's season record repeated on every row ---
$lastRow = 62

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 68
    $ws.Cells.Item($r, 31).Value = 93
    $ws.Cells.Item($r, 32).Value = 0
}
